$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Yasser005@gmail.com"
$ws.Range("D3").Value = "Yasser006@gmail.com"
$ws.Range("D4").Value = "Yasser007@gmail.com"
$ws.Range("D5").Value = "Yasser008@gmail.com"
